$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
$ws.Range('D2').Value = '37.737.99'
$ws.Range('E2').Value = '  -0.21%  '

# Row 3
$ws.Range('D3').Value = '2.078.98'
$ws.Range('E3').Value = '  -0.21%  '

# Row 4
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.17%  '

# Row 5
Set-TextValue $ws.Range('D5') '232.57'
$ws.Range('E5').Value = '  -0.52%  '

# Row 6
$ws.Range('E6').Value = '  -0.31%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
Set-TextValue $ws.Range('D8') '58.10'
$ws.Range('E8').Value = '  -1.67%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.391'
$ws.Range('E9').Value = '  -1.07%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.0781'
$ws.Range('E10').Value = '  -1.13%  '

# Row 11
$ws.Range('E11').Value = '  +0.47%  '

# Row 12
Set-TextValue $ws.Range('D12') '14.87'
$ws.Range('E12').Value = '  +0.45%  '

# Row 13
$ws.Range('D13').Value = '2.385.59'
$ws.Range('E13').Value = '  -0.15%  '

# Row 14
Set-TextValue $ws.Range('D14') '21.23'
$ws.Range('E14').Value = '  -0.28%  '

# Row 15
Set-TextValue $ws.Range('D15') '0.766'
$ws.Range('E15').Value = '  -1.29%  '

# Row 16
Set-TextValue $ws.Range('D16') '5.33'
$ws.Range('E16').Value = '  +0.00%  '

# Row 17
$ws.Range('D17').Value = '2.075.53'
$ws.Range('E17').Value = '  +0.77%  '

# Row 18
$ws.Range('D18').Value = '37.689.94'
$ws.Range('E18').Value = '  -0.11%  '

# Row 19
Set-TextValue $ws.Range('D19') '6.17'
$ws.Range('E19').Value = '  +0.06%  '

# Row 20
Set-TextValue $ws.Range('D20') '70.23'
$ws.Range('E20').Value = '  -2.15%  '

# Row 21
$ws.Range('E21').Value = '  -2.04%  '

# Row 22
Set-TextValue $ws.Range('D22') '227.77'
$ws.Range('E22').Value = '  -0.27%  '

# Row 23
$ws.Range('E23').Value = '  +0.01%  '

# Row 24
Set-TextValue $ws.Range('D24') '2.40'
$ws.Range('E24').Value = '  +0.20%  '

# Row 25
$ws.Range('E25').Value = '  -2.20%  '

# Row 26
Set-TextValue $ws.Range('D26') '9.93'
$ws.Range('E26').Value = '  +3.23%  '

# Row 27
Set-TextValue $ws.Range('D27') '169.50'
$ws.Range('E27').Value = '  -1.05%  '

# Row 28
$ws.Range('E28').Value = '  -3.66%  '

# Row 29
Set-TextValue $ws.Range('D29') '19.37'
$ws.Range('E29').Value = '  -1.00%  '

# Row 30
$ws.Range('E30').Value = '  -2.78%  '

# Row 31
$ws.Range('E31').Value = '  +0.15%  '

# Row 32
$ws.Range('E32').Value = '  -2.97%  '

# Row 33
$ws.Range('E33').Value = '  -0.83%  '

# Row 34
$ws.Range('E34').Value = '  -0.55%  '

# Row 35
$ws.Range('E35').Value = '  +1.15%  '

# Row 36
$ws.Range('E36').Value = '  +0.17%  '

# Row 37
$ws.Range('E37').Value = '  -3.39%  '

# Row 38
$ws.Range('E38').Value = '  -0.05%  '

# Row 39
$ws.Range('E39').Value = '  -1.62%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.0227'
$ws.Range('E40').Value = '  +3.75%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D41') '98.47'
$ws.Range('E41').Value = '  -0.62%  '

# Row 42
$ws.Range('B42').Value = 'Cronos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D42') '0.0960'
$ws.Range('E42').Value = '  -2.31%  '

# Row 43
$ws.Range('E43').Value = '  +0.58%  '

# Row 44
$ws.Range('D44').Value = '1.488.73'

# Row 45
$ws.Range('E45').Value = '  +2.89%  '

# Row 46
Set-TextValue $ws.Range('D46') '16.91'
$ws.Range('E46').Value = '  -2.40%  '

# Row 47
$ws.Range('E47').Value = '  -1.59%  '

# Row 48
$ws.Range('E48').Value = '  -1.89%  '

# Row 49
Set-TextValue $ws.Range('D49') '7.28'
$ws.Range('E49').Value = '  -1.13%  '

# Row 50
$ws.Range('E50').Value = '  -0.83%  '

# Row 51
$ws.Range('D51').Value = '2.269.27'
$ws.Range('E51').Value = '  -0.26%  '

Write-Output "Applied all changes"